# Applies the 2025-09-25 06:34 JST append/update to the "ランサーズ" sheet:
#  - Refresh the "取得日時" timestamp (column A, rows 2-9) from 06:27:39 to 06:34:06
#  - Swap the title (column B) and URL text (column F) between rows 6 and 7,
#    reflecting the two listings exchanging rank/order in the source feed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldTimestamp = "2025-09-25 06:27:39"
$newTimestamp = "2025-09-25 06:34:06"

for ($row = 2; $row -le 9; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    if ($cell.Value2 -eq $oldTimestamp) {
        $cell.Value = $newTimestamp
    }
}

# Rows 6 and 7 swap their title (B) and URL (F) text values.
$b6 = $ws.Cells.Item(6, 2).Value2
$b7 = $ws.Cells.Item(7, 2).Value2
$ws.Cells.Item(6, 2).Value = $b7
$ws.Cells.Item(7, 2).Value = $b6

$f6 = $ws.Cells.Item(6, 6).Value2
$f7 = $ws.Cells.Item(7, 6).Value2
$ws.Cells.Item(6, 6).Value = $f7
$ws.Cells.Item(7, 6).Value = $f6
